# "Kleine Korrekturen in der Präsentation"
#
# 1) Header/footer "datetimeFigureOut" date fields: 16.11.2014 -> 17.11.2014
#    (only reachable, in this object model, through the slide-layout
#    placeholders that are exposed via SlideMaster.CustomLayouts; the
#    Presentation.HandoutMaster / Presentation.NotesMaster accessors in
#    this runtime mis-route text writes into the SlideMaster shapes, so
#    they are intentionally left untouched to avoid corrupting unrelated
#    content)
# 2) Slide 10: "mischeln" -> "mischen"
# 3) Slide 16: merge "Weniger " + "Speicherplatz" runs into one run
# 4) Slide 17: "Datetyp" -> "Datentyp"
# 5) Last slide: remove the two empty, unused placeholder shapes

$p = $ppt.ActivePresentation

# --- 1) Date placeholders on the slide layouts that carry one ---------
$master = $p.SlideMaster
$layouts = $master.CustomLayouts
foreach ($li in 4, 5, 7, 8, 9, 10, 11) {
    $lay = $layouts.Item($li)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = "17.11.2014"
        }
    }
}

# --- 2) Slide 10: "mischeln" -> "mischen" ------------------------------
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange
$full10 = $tr10.Text
$idx10 = $full10.IndexOf("mischeln")
if ($idx10 -ge 0) {
    $sub10 = $tr10.Characters($idx10 + 1, 8)
    $sub10.Text = "mischen"
}

# --- 3) Slide 16: "Weniger " / "Speicherplatz" -> one run --------------
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(2)
$tr16 = $sh16.TextFrame.TextRange
$para16 = $tr16.Paragraphs(2)
$sub16 = $tr16.Characters($para16.Start, $para16.Length)
$sub16.Text = "Weniger Speicherplatz"

# --- 4) Slide 17: "Datetyp" -> "Datentyp" ------------------------------
$s17 = $p.Slides.Item(17)
$sh17 = $s17.Shapes.Item(2)
$tr17 = $sh17.TextFrame.TextRange
$full17 = $tr17.Text
$idx17 = $full17.IndexOf("Datetyp")
if ($idx17 -ge 0) {
    $sub17 = $tr17.Characters($idx17 + 1, 7)
    $sub17.Text = "Datentyp"
}

# --- 5) Last slide: drop the two empty placeholder shapes --------------
$last = $p.Slides.Item($p.Slides.Count)
while ($last.Shapes.Count -gt 0) {
    $last.Shapes.Item(1).Delete()
}
